# Apply cell-value updates per the source diff.
# cryptos.xlsx price/volume refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.490.50"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "2.100.72"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.72"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5225"
$ws.Range("E7").Value = "  -1.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4500"
$ws.Range("E8").Value = "  +2.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.75"
$ws.Range("E9").Value = "  +16.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08929"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.155"
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.47"
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("D13").Value = "2.095.42"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.731"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.709"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.45"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001124"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06621"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.19"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.289"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").Value = "30.544.27"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.33"
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.335"
$ws.Range("E25").Value = "  +3.02%  "
$ws.Range("D26").Value = "2.339.90"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.29"
$ws.Range("E27").Value = "  -2.21%  "
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.31"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.20"
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.201"
$ws.Range("E31").Value = "  +2.88%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.672"
$ws.Range("E33").Value = "  +9.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.165"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.901"
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.45"
$ws.Range("E36").Value = "  +9.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02571"
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06792"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.486"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.71"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6937"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.251"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.02"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6359"
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.266"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.634"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.245"
$ws.Range("E49").Value = "  +7.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.246"
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.10"
$ws.Range("E51").Value = "  -0.52%  "
